$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44729
$ws.Range("J2").Value2 = 300
$ws.Range("K2").Value2 = 14000
$ws.Range("L2").Value2 = 14000
$ws.Range("M2").Value2 = 14000
$ws.Range("P2").Value2 = 778
$ws.Range("D3").Value2 = 44585
$ws.Range("J3").Value2 = 200
$ws.Range("K3").Value2 = 12000
$ws.Range("L3").Value2 = 12000
$ws.Range("M3").Value2 = 12000
$ws.Range("P3").Value2 = 667
$ws.Range("D4").Value2 = 44630
$ws.Range("J4").Value2 = 300
$ws.Range("K4").Value2 = 15000
$ws.Range("L4").Value2 = 15000
$ws.Range("M4").Value2 = 15000
$ws.Range("P4").Value2 = 833
$ws.Range("D5").Value2 = 44243
$ws.Range("I5").Value = "Especial"
$ws.Range("K5").Value2 = 12000
$ws.Range("L5").Value2 = 12000
$ws.Range("M5").Value2 = 12000
$ws.Range("P5").Value2 = 667
$ws.Range("D6").Value2 = 44243
$ws.Range("J6").Value2 = 300
$ws.Range("K6").Value2 = 10000
$ws.Range("L6").Value2 = 10000
$ws.Range("M6").Value2 = 10000
$ws.Range("P6").Value2 = 556
$ws.Range("D7").Value2 = 44243
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value2 = 150
$ws.Range("K7").Value2 = 8000
$ws.Range("L7").Value2 = 8000
$ws.Range("M7").Value2 = 8000
$ws.Range("P7").Value2 = 444
$ws.Range("D8").Value2 = 44750
$ws.Range("I8").Value = "Primera"
$ws.Range("K8").Value2 = 15000
$ws.Range("L8").Value2 = 15000
$ws.Range("M8").Value2 = 15000
$ws.Range("P8").Value2 = 833
$ws.Range("D10").Value2 = 44742
$ws.Range("J10").Value2 = 300
$ws.Range("K10").Value2 = 14000
$ws.Range("L10").Value2 = 14000
$ws.Range("M10").Value2 = 14000
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value2 = 778
$ws.Range("D11").Value2 = 44719
$ws.Range("D12").Value2 = 44631
$ws.Range("J12").Value2 = 300
$ws.Range("K12").Value2 = 15000
$ws.Range("L12").Value2 = 15000
$ws.Range("M12").Value2 = 15000
$ws.Range("P12").Value2 = 833
$ws.Range("D13").Value2 = 44245
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value2 = 300
$ws.Range("K13").Value2 = 12000
$ws.Range("L13").Value2 = 12000
$ws.Range("M13").Value2 = 12000
$ws.Range("P13").Value2 = 667
$ws.Range("D14").Value2 = 44245
$ws.Range("I14").Value = "Segunda"
$ws.Range("J14").Value2 = 200
$ws.Range("K14").Value2 = 10000
$ws.Range("L14").Value2 = 10000
$ws.Range("M14").Value2 = 10000
$ws.Range("P14").Value2 = 556
$ws.Range("D15").Value2 = 44635
$ws.Range("D16").Value2 = 44383
$ws.Range("J16").Value2 = 300
$ws.Range("K16").Value2 = 16000
$ws.Range("L16").Value2 = 16000
$ws.Range("M16").Value2 = 16000
$ws.Range("P16").Value2 = 889
$ws.Range("D17").Value2 = 44383
$ws.Range("K17").Value2 = 12000
$ws.Range("L17").Value2 = 12000
$ws.Range("M17").Value2 = 12000
$ws.Range("P17").Value2 = 667
$ws.Range("D18").Value2 = 44249
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value2 = 400
$ws.Range("K18").Value2 = 12000
$ws.Range("L18").Value2 = 12000
$ws.Range("M18").Value2 = 12000
$ws.Range("P18").Value2 = 667
$ws.Range("D19").Value2 = 44249
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value2 = 200
$ws.Range("K19").Value2 = 10000
$ws.Range("L19").Value2 = 10000
$ws.Range("M19").Value2 = 10000
$ws.Range("P19").Value2 = 556
$ws.Range("D20").Value2 = 44396
$ws.Range("J20").Value2 = 250
$ws.Range("K20").Value2 = 15000
$ws.Range("L20").Value2 = 15000
$ws.Range("M20").Value2 = 15000
$ws.Range("P20").Value2 = 833
$ws.Range("D21").Value2 = 44396
$ws.Range("K21").Value2 = 12000
$ws.Range("L21").Value2 = 12000
$ws.Range("M21").Value2 = 12000
$ws.Range("P21").Value2 = 667
$ws.Range("D22").Value2 = 44721
$ws.Range("K22").Value2 = 15000
$ws.Range("L22").Value2 = 15000
$ws.Range("M22").Value2 = 15000
$ws.Range("P22").Value2 = 833
$ws.Range("D23").Value2 = 44627
$ws.Range("J23").Value2 = 300
$ws.Range("D24").Value2 = 44753
$ws.Range("J24").Value2 = 300
$ws.Range("D25").Value2 = 44614
$ws.Range("N25").Value = "$/caja 18 kilos granel"
$ws.Range("D26").Value2 = 44754
$ws.Range("J26").Value2 = 400
$ws.Range("K26").Value2 = 15000
$ws.Range("L26").Value2 = 15000
$ws.Range("M26").Value2 = 15000
$ws.Range("P26").Value2 = 833
$ws.Range("D27").Value2 = 44725
$ws.Range("J27").Value2 = 400
$ws.Range("O27").Value = "Región de Arica y Parinacota"
$ws.Range("D28").Value2 = 44736
$ws.Range("J28").Value2 = 200
$ws.Range("D29").Value2 = 44746
$ws.Range("J29").Value2 = 500
$ws.Range("N29").Value = "$/bandeja 18 kilos"
$ws.Range("D31").Value2 = 44235
$ws.Range("J31").Value2 = 400
$ws.Range("K31").Value2 = 13000
$ws.Range("L31").Value2 = 13000
$ws.Range("M31").Value2 = 13000
$ws.Range("P31").Value2 = 722
$ws.Range("D32").Value2 = 44235
$ws.Range("I32").Value = "Segunda"
$ws.Range("J32").Value2 = 200
$ws.Range("K32").Value2 = 11000
$ws.Range("L32").Value2 = 11000
$ws.Range("M32").Value2 = 11000
$ws.Range("P32").Value2 = 611
$ws.Range("D33").Value2 = 44235
$ws.Range("I33").Value = "Tercera"
$ws.Range("J33").Value2 = 100
$ws.Range("K33").Value2 = 9000
$ws.Range("L33").Value2 = 9000
$ws.Range("M33").Value2 = 9000
$ws.Range("P33").Value2 = 500
$ws.Range("D34").Value2 = 44726
$ws.Range("K34").Value2 = 14000
$ws.Range("L34").Value2 = 14000
$ws.Range("M34").Value2 = 14000
$ws.Range("P34").Value2 = 778
$ws.Range("D35").Value2 = 44628
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value2 = 300
$ws.Range("K35").Value2 = 15000
$ws.Range("L35").Value2 = 15000
$ws.Range("M35").Value2 = 15000
$ws.Range("P35").Value2 = 833
$ws.Range("D37").Value2 = 44238
$ws.Range("K37").Value2 = 12000
$ws.Range("L37").Value2 = 12000
$ws.Range("M37").Value2 = 12000
$ws.Range("P37").Value2 = 667
$ws.Range("D38").Value2 = 44238
$ws.Range("I38").Value = "Segunda"
$ws.Range("J38").Value2 = 200
$ws.Range("K38").Value2 = 10000
$ws.Range("L38").Value2 = 10000
$ws.Range("M38").Value2 = 10000
$ws.Range("P38").Value2 = 556
$ws.Range("I39").Value = "Tercera"
$ws.Range("J39").Value2 = 50
$ws.Range("K39").Value2 = 8000
$ws.Range("L39").Value2 = 8000
$ws.Range("M39").Value2 = 8000
$ws.Range("P39").Value2 = 444
$ws.Range("D40").Value2 = 44596
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value2 = 150
$ws.Range("K40").Value2 = 14000
$ws.Range("L40").Value2 = 14000
$ws.Range("M40").Value2 = 14000
$ws.Range("P40").Value2 = 778
$ws.Range("D41").Value2 = 44757
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value2 = 300
$ws.Range("K41").Value2 = 15000
$ws.Range("L41").Value2 = 15000
$ws.Range("M41").Value2 = 15000
$ws.Range("P41").Value2 = 833
$ws.Range("D42").Value2 = 44391
$ws.Range("J42").Value2 = 400
